# Weekly fruit/hortaliza update:
#  - a new observation is inserted as the new first data row (row 18),
#    pushing the existing row 18..49 block down by one row
#  - a second new observation is inserted a bit further down (row 43 in
#    the resulting layout), pushing the remaining rows down by one more
# Net effect: 49 data rows -> 51 data rows (2 new rows), dimension A1:R50 -> A1:R52

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at row 18 ---------------------------------------
$ws.Rows("18:18").Insert()

$row = 18
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44645
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112043
$ws.Cells.Item($row, 7).Value = "Pepino dulce"
$ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 80
$ws.Cells.Item($row, 11).Value = 18000
$ws.Cells.Item($row, 12).Value = 18000
$ws.Cells.Item($row, 13).Value = 18000
$ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 1000
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"

# --- Insert second new row at row 43 (post first insert) -------------------
$ws.Rows("43:43").Insert()

$row = 43
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44644
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112043
$ws.Cells.Item($row, 7).Value = "Pepino dulce"
$ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 9).Value = "Especial"
$ws.Cells.Item($row, 10).Value = 50
$ws.Cells.Item($row, 11).Value = 21000
$ws.Cells.Item($row, 12).Value = 21000
$ws.Cells.Item($row, 13).Value = 21000
$ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 1167
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"
